$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume(1h) columns (D/E) with the latest scrape.
# Plain decimal-looking prices are written with a leading apostrophe so Excel
# keeps them as literal text (matching the sheet's existing text-cell layout)
# instead of silently coercing them to floating-point numbers.
$ws.Range("D2").Value = '27.991.97'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '1.647.67'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''213.71'
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").Value = '''0.528'
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''23.45'
$ws.Range("E8").Value = '  +2.78%  '
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("D11").Value = '''0.0872'
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").Value = '1.881.51'
$ws.Range("D13").Value = '1.642.63'
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '''0.564'
$ws.Range("E15").Value = '  +2.78%  '
$ws.Range("D16").Value = '''65.62'
$ws.Range("E16").Value = '  +0.83%  '
$ws.Range("D17").Value = '27.998.99'
$ws.Range("D18").Value = '''232.92'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '''7.71'
$ws.Range("E19").Value = '  +2.71%  '
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("E22").Value = '  +4.92%  '
$ws.Range("D23").Value = '''4.39'
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("D24").Value = '''2.17'
$ws.Range("E24").Value = '  +4.65%  '
$ws.Range("D25").Value = '''152.62'
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("D26").Value = '''6.93'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("E32").Value = '  +3.03%  '
$ws.Range("D33").Value = '1.449.75'
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").Value = '''3.09'
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("E37").Value = '  +3.37%  '
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").Value = '''0.922'
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").Value = '''69.52'
$ws.Range("E41").Value = '  +2.81%  '
$ws.Range("E42").Value = '  +3.27%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  +0.92%  '
$ws.Range("D45").Value = '''2.22'
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("D46").Value = '''5.40'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  +5.08%  '
$ws.Range("D48").Value = '1.789.88'
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("D49").Value = '''88.96'
$ws.Range("E49").Value = '  +2.87%  '
$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("E51").Value = '  +0.43%  '
